# The "date" column (A) currently stores the 1st-of-month of a quarter
# start (Jan/Apr/Jul/Oct). The fix re-indexes each quarter stamp to the
# 15th of the following month (Feb/May/Aug/Nov), keeping the same year.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$epoch = Get-Date -Year 1899 -Month 12 -Day 30 -Hour 0 -Minute 0 -Second 0

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $old = $cell.Value2
    if ($old -eq $null -or $old -eq "") { continue }

    $d = $epoch.AddDays($old)
    $newDate = Get-Date -Year $d.Year -Month ($d.Month + 1) -Day 15 -Hour 0 -Minute 0 -Second 0
    $newSerial = [int]$newDate.ToOADate()
    $cell.Value = $newSerial
}
